# Update the "想去人数" (F column) counts on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - row => new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 320
$ws1.Range("F6").Value  = 664
$ws1.Range("F7").Value  = 269
$ws1.Range("F12").Value = 3366
$ws1.Range("F20").Value = 670
$ws1.Range("F23").Value = 54
$ws1.Range("F24").Value = 47
$ws1.Range("F26").Value = 2389
$ws1.Range("F27").Value = 4931
$ws1.Range("F31").Value = 1266
$ws1.Range("F32").Value = 269
$ws1.Range("F33").Value = 2188
$ws1.Range("F37").Value = 75
$ws1.Range("F39").Value = 305
$ws1.Range("F40").Value = 451

# Sheet "全部类型" (sheet4) - row => new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 320
$ws4.Range("F6").Value  = 664
$ws4.Range("F7").Value  = 269
$ws4.Range("F12").Value = 3366
$ws4.Range("F21").Value = 670
$ws4.Range("F24").Value = 54
$ws4.Range("F25").Value = 47
$ws4.Range("F27").Value = 2389
$ws4.Range("F28").Value = 4931
$ws4.Range("F32").Value = 1266
$ws4.Range("F33").Value = 269
$ws4.Range("F34").Value = 2188
$ws4.Range("F38").Value = 75
$ws4.Range("F40").Value = 305
$ws4.Range("F41").Value = 451
